# I2C-SMBus_Diagram.pptx edits
#  - retitle the diagram shape: "I2C/SMBus Diagram" -> "InternalCharger: I2C/SMBus Interfaces"
#    (and widen/shift it to fit the new, longer title)
#  - swap the "Cell A" / "Cell B" labels on the battery-cell shapes
#    (revert of an off-by-one fix on the TWI master/slave modules)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title shape: "Rectangle 1" (shape id 2) -----------------------------
$title = $s.Shapes.Item(3)

# Reposition/resize (EMU -> points, 914400 EMU per inch, 72 pt per inch)
$title.Left   = 2681515 / 914400 * 72
$title.Top    = 311287  / 914400 * 72
$title.Width  = 3775842 / 914400 * 72
$title.Height = 369332  / 914400 * 72

$tr = $title.TextFrame.TextRange

# Prepend a new run: "InternalCharger"
[void]$tr.InsertBefore("InternalCharger")

# "I2C/" run becomes ": I2C/" (now starts right after "InternalCharger")
$tr.Characters(("InternalCharger").Length + 1, 4).Text = ": I2C/"

# " Diagram" run (after "InternalCharger" + ": I2C/" + "SMBus") becomes " Interfaces"
$prefixLen = ("InternalCharger" + ": I2C/" + "SMBus").Length
$tr.Characters($prefixLen + 1, 8).Text = " Interfaces"

# --- Cell A / Cell B swap -------------------------------------------------
$cellA = $s.Shapes.Item(23)   # "Rounded Rectangle 105" -> currently "CellA"
$cellB = $s.Shapes.Item(24)   # "Rounded Rectangle 112" -> currently "CellB"

$cellA.TextFrame.TextRange.Characters($cellA.TextFrame.TextRange.Length, 1).Text = "B"
$cellB.TextFrame.TextRange.Characters($cellB.TextFrame.TextRange.Length, 1).Text = "A"
